# Added headless flag ("execute" column, Y/N) to the Screens sheet and
# moved the active selection to E4 (optimized Menu search flow).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Screens")

# Populate the new "execute" (headless) flag column with Y/N values.
$ws.Range("E2").Value  = "N"
$ws.Range("E3").Value  = "N"
$ws.Range("E4").Value  = "Y"
$ws.Range("E5").Value  = "Y"
$ws.Range("E6").Value  = "Y"
$ws.Range("E7").Value  = "N"
$ws.Range("E8").Value  = "N"
$ws.Range("E9").Value  = "N"
$ws.Range("E10").Value = "N"
$ws.Range("E11").Value = "N"
$ws.Range("E12").Value = "N"

# Update the active selection/cell as saved in the workbook view.
$ws.Range("E4").Select() | Out-Null
